# Applies the cryptos.xlsx data refresh described in the commit:
# "Updated cryptos list on Wed Jan  3 19:45:38 UTC 2024 with GitHub Actions"
#
# Column D ("Price") values look like numbers (including some using
# "thousands." separators, e.g. 43.003.68) but must stay as plain TEXT,
# exactly like the original workbook (t="inlineStr"/shared string, not t="n").
# Writing a plain numeric-looking string via .Value lets Excel silently
# coerce it into a real number, so we prefix it with a leading apostrophe
# (forces text entry) and then reset the cell .Style to "Normal" so no stray
# number-format / quote-prefix styling is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'43.003.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.74%  "

# Row 3
$ws.Range("D3").Value = "'2.218.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.13%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "'316.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.00%  "

# Row 6
$ws.Range("D6").Value = "'98.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.75%  "

# Row 7
$ws.Range("D7").Value = "'0.583"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.64%  "

# Row 8
$ws.Range("E8").Value = "  +0.12%  "

# Row 9
$ws.Range("D9").Value = "'0.562"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.62%  "

# Row 10
$ws.Range("D10").Value = "'36.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -10.30%  "

# Row 11
$ws.Range("D11").Value = "'54.30"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.06%  "

# Row 12
$ws.Range("D12").Value = "'0.0827"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.73%  "

# Row 13
$ws.Range("D13").Value = "'7.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.77%  "

# Row 14
$ws.Range("D14").Value = "'0.106"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.49%  "

# Row 15
$ws.Range("D15").Value = "'0.863"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -11.41%  "

# Row 16
$ws.Range("D16").Value = "'2.558.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.18%  "

# Row 17
$ws.Range("D17").Value = "'14.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.75%  "

# Row 18
$ws.Range("D18").Value = "'2.214.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.30%  "

# Row 19
$ws.Range("D19").Value = "'42.856.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.00%  "

# Row 20
$ws.Range("D20").Value = "'14.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.06%  "

# Row 21
$ws.Range("E21").Value = "  -8.97%  "

# Row 22
$ws.Range("D22").Value = "'6.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -11.78%  "

# Row 23
$ws.Range("E23").Value = "  -10.67%  "

# Row 24
$ws.Range("E24").Value = "  -9.33%  "

# Row 25
$ws.Range("D25").Value = "'237.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.41%  "

# Row 26
$ws.Range("D26").Value = "'2.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.86%  "

# Row 27
$ws.Range("E27").Value = "  +0.12%  "

# Row 28
$ws.Range("D28").Value = "'10.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.05%  "

# Row 29
$ws.Range("D29").Value = "'2.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.69%  "

# Row 30
$ws.Range("D30").Value = "'6.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -13.04%  "

# Row 31
$ws.Range("D31").Value = "'20.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.20%  "

# Row 32
$ws.Range("D32").Value = "'0.0871"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.96%  "

# Row 33
$ws.Range("D33").Value = "'33.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.86%  "

# Row 34
$ws.Range("D34").Value = "'155.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.74%  "

# Row 35
$ws.Range("D35").Value = "'2.76"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.47%  "

# Row 36
$ws.Range("D36").Value = "'3.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.22%  "

# Row 37
$ws.Range("E37").Value = "  +15.33%  "

# Row 38
$ws.Range("E38").Value = "  -6.30%  "

# Row 39
$ws.Range("E39").Value = "  -5.59%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.102"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.85%  "

# Row 41
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").Value = "'3.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.08%  "

# Row 42
$ws.Range("D42").Value = "'0.0326"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.62%  "

# Row 43
$ws.Range("D43").Value = "'1.888.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.36%  "

# Row 44
$ws.Range("E44").Value = "  +0.15%  "

# Row 45
$ws.Range("D45").Value = "'90.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.30%  "

# Row 46
$ws.Range("D46").Value = "'12.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.81%  "

# Row 47
$ws.Range("E47").Value = "  -9.52%  "

# Row 48
$ws.Range("D48").Value = "'5.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.02%  "

# Row 49
$ws.Range("D49").Value = "'60.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -12.40%  "

# Row 50
$ws.Range("D50").Value = "'75.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.21%  "

# Row 51
$ws.Range("B51").Value = "SEI"
$ws.Range("C51").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D51").Value = "'0.865"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +15.54%  "
